# Auto-generated Excel COM-interop script applying value changes
# described by the commit diff across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 100008
$ws.Range("I21").Value = 100008
$ws.Range("K21").Value = 100008
$ws.Range("M21").Value = -99540
$ws.Range("H23").Value = 100008
$ws.Range("I23").Value = 100008
$ws.Range("K23").Value = 100008
$ws.Range("M23").Value = -99774
$ws.Range("H113").Value = 8629.4375
$ws.Range("I113").Value = 11953
$ws.Range("J113").Value = 3090.1667
$ws.Range("K113").Value = 11953
$ws.Range("L113").Value = 3090.1667
$ws.Range("M113").Value = -8699
$ws.Range("N113").Value = -9598.1667
$ws.Range("H127").Value = 860.9474
$ws.Range("I127").Value = 709.9375
$ws.Range("K127").Value = 2129.8125
$ws.Range("M127").Value = 2830.1875
$ws.Range("H137").Value = 70890.914
$ws.Range("I137").Value = 70890.914
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 212672.742
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -210122.742
$ws.Range("N137").ClearContents() | Out-Null
$ws.Range("H138").Value = 3910.121
$ws.Range("I138").Value = 1972.9166
$ws.Range("J138").Value = 4177.322
$ws.Range("K138").Value = 5918.7498
$ws.Range("L138").Value = 12531.966
$ws.Range("M138").Value = -778.7497999999996
$ws.Range("N138").Value = -22811.966
$ws.Range("H141").Value = 1042.5
$ws.Range("I141").Value = 1042.5
$ws.Range("K141").Value = 3127.5
$ws.Range("M141").Value = 2052.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 430
$ws.Range("I12").Value = 430
$ws.Range("K12").Value = 430
$ws.Range("M12").Value = -257
$ws.Range("H61").Value = 2622
$ws.Range("I61").Value = 1967.8636
$ws.Range("K61").Value = 1967.8636
$ws.Range("M61").Value = -1755.8636
$ws.Range("H102").Value = 1638.5625
$ws.Range("I102").Value = 789.75
$ws.Range("K102").Value = 789.75
$ws.Range("M102").Value = 832.25
$ws.Range("H122").Value = 5436.8096
$ws.Range("I122").Value = 3935.3076
$ws.Range("J122").Value = 7876.75
$ws.Range("K122").Value = 11805.9228
$ws.Range("L122").Value = 23630.25
$ws.Range("M122").Value = -9355.9228
$ws.Range("N122").Value = -28530.25
$ws.Range("H129").Value = 129932.664
$ws.Range("J129").Value = 129932.664
$ws.Range("L129").Value = 129932.664
$ws.Range("N129").Value = -139932.664
$ws.Range("H132").Value = 2537.4443
$ws.Range("I132").Value = 2175.6052
$ws.Range("K132").Value = 6526.8156
$ws.Range("M132").Value = -3996.8156
$ws.Range("H136").Value = 2622
$ws.Range("I136").Value = 1967.8636
$ws.Range("K136").Value = 5903.5908
$ws.Range("M136").Value = -3353.5908
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1148.8334
$ws.Range("I94").Value = 925.61536
$ws.Range("K94").Value = 925.61536
$ws.Range("M94").Value = -474.61536
$ws.Range("H105").Value = 2050.7727
$ws.Range("I105").Value = 1905.9
$ws.Range("K105").Value = 1905.9
$ws.Range("M105").Value = -158.9000000000001
$ws.Range("H134").Value = 5137.364
$ws.Range("I134").Value = 3128.1428
$ws.Range("K134").Value = 9384.428400000001
$ws.Range("M134").Value = -6849.428400000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 40448.855
$ws.Range("J28").Value = 40448.855
$ws.Range("L28").Value = 40448.855
$ws.Range("N28").Value = -40938.855
$ws.Range("H31").Value = 5440.2896
$ws.Range("I31").Value = 1866
$ws.Range("J31").Value = 8039.773
$ws.Range("K31").Value = 1866
$ws.Range("L31").Value = 8039.773
$ws.Range("M31").Value = -1571
$ws.Range("N31").Value = -8629.773000000001
$ws.Range("H34").Value = 5440.2896
$ws.Range("I34").Value = 1866
$ws.Range("J34").Value = 8039.773
$ws.Range("K34").Value = 1866
$ws.Range("L34").Value = 8039.773
$ws.Range("M34").Value = -1664
$ws.Range("N34").Value = -8443.773000000001
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents() | Out-Null
$ws.Range("H58").Value = 3549.7273
$ws.Range("I58").Value = 3528.2856
$ws.Range("K58").Value = 3528.2856
$ws.Range("M58").Value = -3325.2856
$ws.Range("H132").Value = 3381.6072
$ws.Range("I132").Value = 3520.7083
$ws.Range("K132").Value = 10562.1249
$ws.Range("M132").Value = -8032.124899999999
$ws.Range("H134").Value = 2017.1613
$ws.Range("I134").Value = 1701.2593
$ws.Range("K134").Value = 5103.7779
$ws.Range("M134").Value = -2568.7779
$ws.Range("H136").Value = 3549.7273
$ws.Range("I136").Value = 3528.2856
$ws.Range("K136").Value = 10584.8568
$ws.Range("M136").Value = -8034.856800000001
$ws.Range("H137").Value = 90000.664
$ws.Range("I137").Value = 30000
$ws.Range("K137").Value = 30000
$ws.Range("M137").Value = -24900
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 722.6
$ws.Range("I47").Value = 704.5
$ws.Range("J47").Value = 795
$ws.Range("K47").Value = 2113.5
$ws.Range("L47").Value = 2385
$ws.Range("M47").Value = -1682.5
$ws.Range("N47").Value = -3247
$ws.Range("H68").Value = 1458
$ws.Range("I68").Value = 824.5
$ws.Range("J68").Value = 1639
$ws.Range("K68").Value = 2473.5
$ws.Range("L68").Value = 4917
$ws.Range("M68").Value = -1662.5
$ws.Range("N68").Value = -6539
$ws.Range("H71").Value = 1458
$ws.Range("I71").Value = 824.5
$ws.Range("J71").Value = 1639
$ws.Range("K71").Value = 7420.5
$ws.Range("L71").Value = 14751
$ws.Range("M71").Value = -3364.5
$ws.Range("N71").Value = -22863
$ws.Range("H82").Value = 2000
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5594
$ws.Range("H85").Value = 2000
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4596
$ws.Range("H109").Value = 2047.2222
$ws.Range("I109").Value = 1070.8334
$ws.Range("K109").Value = 3212.5002
$ws.Range("M109").Value = -2172.5002
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents() | Out-Null
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents() | Out-Null
$ws.Range("H124").Value = 4300
$ws.Range("J124").Value = 5000
$ws.Range("L124").Value = 15000
$ws.Range("N124").Value = -24820
$ws.Range("H129").Value = 1875
$ws.Range("I129").Value = 737.2
$ws.Range("K129").Value = 2211.6
$ws.Range("M129").Value = 2788.4
$ws.Range("H130").Value = 2976.8572
$ws.Range("I130").Value = 2167.6
$ws.Range("K130").Value = 6502.799999999999
$ws.Range("M130").Value = -1482.799999999999
$ws.Range("H131").Value = 1316.7142
$ws.Range("I131").Value = 798.1429000000001
$ws.Range("J131").Value = 1835.2858
$ws.Range("K131").Value = 2394.4287
$ws.Range("L131").Value = 5505.857400000001
$ws.Range("M131").Value = 2645.5713
$ws.Range("N131").Value = -15585.8574
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H106").Value = 23500
$ws.Range("J106").Value = 23500
$ws.Range("L106").Value = 23500
$ws.Range("N106").Value = -26024
$ws.Range("H113").Value = 1561.6
$ws.Range("I113").Value = 1398.3334
$ws.Range("J113").Value = 1806.5
$ws.Range("K113").Value = 1398.3334
$ws.Range("L113").Value = 1806.5
$ws.Range("M113").Value = 771.6666
$ws.Range("N113").Value = -6146.5
$ws.Range("H132").Value = 2877.762
$ws.Range("I132").Value = 2446.9697
$ws.Range("J132").Value = 4457.3335
$ws.Range("K132").Value = 7340.909100000001
$ws.Range("L132").Value = 13372.0005
$ws.Range("M132").Value = -4810.909100000001
$ws.Range("N132").Value = -18432.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5415.55
$ws.Range("I132").Value = 5477.2188
$ws.Range("J132").Value = 5168.875
$ws.Range("K132").Value = 16431.6564
$ws.Range("L132").Value = 15506.625
$ws.Range("M132").Value = -13901.6564
$ws.Range("N132").Value = -20566.625
$ws.Range("H136").Value = 2349.7856
$ws.Range("I136").Value = 1611.625
$ws.Range("K136").Value = 4834.875
$ws.Range("M136").Value = -2284.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4047.8235
$ws.Range("I126").Value = 4220.8667
$ws.Range("K126").Value = 12662.6001
$ws.Range("M126").Value = -10192.6001

Write-Host "Applied all cell updates."
